$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.719.92'
$ws.Cells.Item(2, 5).Value = '  -0.51%  '

$ws.Cells.Item(3, 4).Value = '2.440.27'
$ws.Cells.Item(3, 5).Value = '  -1.07%  '

$ws.Cells.Item(4, 5).Value = '  +0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '507.87'
$ws.Cells.Item(5, 5).Value = '  -1.73%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '129.30'
$ws.Cells.Item(6, 5).Value = '  -1.95%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.20%  '

$ws.Cells.Item(8, 5).Value = '  -1.18%  '

$ws.Cells.Item(9, 4).Value = '2.456.42'
$ws.Cells.Item(9, 5).Value = '  -0.53%  '

$ws.Cells.Item(10, 5).Value = '  -0.06%  '

$ws.Cells.Item(11, 5).Value = '  -3.76%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '5.18'
$ws.Cells.Item(12, 5).Value = '  -4.17%  '

$ws.Cells.Item(13, 5).Value = '  -3.05%  '

$ws.Cells.Item(14, 4).Value = '2.872.21'
$ws.Cells.Item(14, 5).Value = '  -1.14%  '

$ws.Cells.Item(15, 4).Value = '57.652.73'
$ws.Cells.Item(15, 5).Value = '  -0.48%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '21.95'
$ws.Cells.Item(16, 5).Value = '  -0.41%  '

$ws.Cells.Item(17, 5).Value = '  -2.50%  '

$ws.Cells.Item(18, 4).Value = '2.445.81'
$ws.Cells.Item(18, 5).Value = '  -1.00%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '10.49'
$ws.Cells.Item(19, 5).Value = '  -3.25%  '

$ws.Cells.Item(20, 5).Value = '  -1.09%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '315.31'
$ws.Cells.Item(21, 5).Value = '  -1.10%  '

$ws.Cells.Item(22, 5).Value = '  -0.02%  '

$ws.Cells.Item(23, 5).Value = '  -1.30%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '63.39'
$ws.Cells.Item(24, 5).Value = '  -1.23%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.407'
$ws.Cells.Item(25, 5).Value = '  -0.31%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.995'
$ws.Cells.Item(26, 5).Value = '  -0.51%  '

$ws.Cells.Item(27, 5).Value = '  -1.25%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.28'
$ws.Cells.Item(28, 5).Value = '  -1.48%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '169.89'
$ws.Cells.Item(29, 5).Value = '  +3.14%  '

$ws.Cells.Item(30, 4).Value = '0.0₃0725'
$ws.Cells.Item(30, 5).Value = '  -2.79%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.27'
$ws.Cells.Item(31, 5).Value = '  -1.60%  '

$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.66'
$ws.Cells.Item(32, 5).Value = '  -2.27%  '

$ws.Cells.Item(33, 2).Value = 'Fetch.AI'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.16'
$ws.Cells.Item(33, 5).Value = '  +2.72%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.998'
$ws.Cells.Item(34, 5).Value = '  +0.00%  '

$ws.Cells.Item(35, 5).Value = '  -0.19%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '17.70'
$ws.Cells.Item(36, 5).Value = '  -2.31%  '

$ws.Cells.Item(37, 5).Value = '  -4.63%  '

$ws.Cells.Item(38, 5).Value = '  -0.33%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '36.30'
$ws.Cells.Item(39, 5).Value = '  -0.44%  '

$ws.Cells.Item(40, 5).Value = '  -1.65%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.770'
$ws.Cells.Item(41, 5).Value = '  -2.42%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '273.20'
$ws.Cells.Item(42, 5).Value = '  -0.46%  '

$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '3.39'
$ws.Cells.Item(43, 5).Value = '  -2.45%  '

$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '5.01'
$ws.Cells.Item(44, 5).Value = '  +1.41%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.581'
$ws.Cells.Item(45, 5).Value = '  -1.49%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.0908'
$ws.Cells.Item(46, 5).Value = '  -0.07%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '120.21'
$ws.Cells.Item(47, 5).Value = '  -5.26%  '

$ws.Cells.Item(48, 5).Value = '  -0.97%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '17.19'
$ws.Cells.Item(49, 5).Value = '  -3.19%  '

$ws.Cells.Item(50, 5).Value = '  -2.27%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '16.69'
$ws.Cells.Item(51, 5).Value = '  -2.02%  '
